# Update column F (dSF) values after re-pulling data / recalculating mean
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -5
$ws.Range("F3").Value  = 2
$ws.Range("F4").Value  = 6
$ws.Range("F5").Value  = -2
$ws.Range("F6").Value  = -5
$ws.Range("F7").Value  = 6
$ws.Range("F8").Value  = -1
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = 1
